$d = $word.ActiveDocument

# Remove the first paragraph ("!!NOT FINISHED!!") entirely, including its
# paragraph mark, so the document begins with the next paragraph.
$firstPara = $d.Paragraphs(1)
$firstPara.Range.Delete()
